# Auto-generated edit script applying the cryptos.xlsx data refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'68.025.92"
$ws.Range("E2").Value = "'  +0.28%  "

# Row 3
$ws.Range("D3").Value = "'3.244.84"
$ws.Range("E3").Value = "'  -0.05%  "

# Row 4
$ws.Range("E4").Value = "'  +0.00%  "

# Row 5
$ws.Range("D5").Value = "'581.81"
$ws.Range("E5").Value = "'  +0.55%  "

# Row 6
$ws.Range("D6").Value = "'184.25"
$ws.Range("E6").Value = "'  +0.76%  "

# Row 7
$ws.Range("E7").Value = "'  +0.01%  "

# Row 8
$ws.Range("E8").Value = "'  +0.90%  "

# Row 9
$ws.Range("E9").Value = "'  -3.30%  "

# Row 10
$ws.Range("E10").Value = "'  -1.03%  "

# Row 11
$ws.Range("D11").Value = "'0.416"
$ws.Range("E11").Value = "'  +0.31%  "

# Row 12
$ws.Range("D12").Value = "'3.808.01"
$ws.Range("E12").Value = "'  -0.03%  "

# Row 13
$ws.Range("D13").Value = "'0.138"
$ws.Range("E13").Value = "'  +0.18%  "

# Row 14
$ws.Range("D14").Value = "'27.87"
$ws.Range("E14").Value = "'  -2.68%  "

# Row 15
$ws.Range("D15").Value = "'68.041.55"
$ws.Range("E15").Value = "'  +0.29%  "

# Row 16
$ws.Range("E16").Value = "'  -1.06%  "

# Row 17
$ws.Range("D17").Value = "'3.219.16"
$ws.Range("E17").Value = "'  -1.00%  "

# Row 18
$ws.Range("E18").Value = "'  -0.40%  "

# Row 19
$ws.Range("D19").Value = "'13.47"
$ws.Range("E19").Value = "'  -0.65%  "

# Row 20
$ws.Range("D20").Value = "'396.40"
$ws.Range("E20").Value = "'  +4.52%  "

# Row 21
$ws.Range("D21").Value = "'7.60"
$ws.Range("E21").Value = "'  -0.50%  "

# Row 22
$ws.Range("E22").Value = "'  +0.12%  "

# Row 23
$ws.Range("D23").Value = "'71.33"
$ws.Range("E23").Value = "'  -0.08%  "

# Row 24
$ws.Range("D24").Value = "'0.516"
$ws.Range("E24").Value = "'  +0.61%  "

# Row 25
$ws.Range("E25").Value = "'  -0.73%  "

# Row 26
$ws.Range("D26").Value = "'0.186"
$ws.Range("E26").Value = "'  +2.49%  "

# Row 27
$ws.Range("D27").Value = "'9.62"
$ws.Range("E27").Value = "'  -3.13%  "

# Row 28
$ws.Range("E28").Value = "'  -0.01%  "

# Row 29
$ws.Range("E29").Value = "'  -0.31%  "

# Row 30
$ws.Range("D30").Value = "'5.60"
$ws.Range("E30").Value = "'  -1.25%  "

# Row 31
$ws.Range("D31").Value = "'22.81"
$ws.Range("E31").Value = "'  -0.11%  "

# Row 32
$ws.Range("D32").Value = "'7.01"
$ws.Range("E32").Value = "'  +0.03%  "

# Row 33
$ws.Range("D33").Value = "'1.26"
$ws.Range("E33").Value = "'  +0.03%  "

# Row 34
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "'  +0.03%  "

# Row 35
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "'161.75"
$ws.Range("E35").Value = "'  -0.37%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.49"
$ws.Range("E36").Value = "'  -5.57%  "

# Row 37
$ws.Range("D37").Value = "'1.90"
$ws.Range("E37").Value = "'  +2.97%  "

# Row 38
$ws.Range("D38").Value = "'0.813"
$ws.Range("E38").Value = "'  -2.99%  "

# Row 39
$ws.Range("D39").Value = "'26.55"
$ws.Range("E39").Value = "'  +0.33%  "

# Row 40
$ws.Range("D40").Value = "'4.60"
$ws.Range("E40").Value = "'  +0.29%  "

# Row 41
$ws.Range("E41").Value = "'  -2.91%  "

# Row 42
$ws.Range("E42").Value = "'  -4.05%  "

# Row 43
$ws.Range("D43").Value = "'41.18"
$ws.Range("E43").Value = "'  +0.01%  "

# Row 44
$ws.Range("D44").Value = "'0.0683"
$ws.Range("E44").Value = "'  -0.40%  "

# Row 45
$ws.Range("D45").Value = "'25.03"
$ws.Range("E45").Value = "'  -1.74%  "

# Row 46
$ws.Range("D46").Value = "'2.606.15"
$ws.Range("E46").Value = "'  -0.95%  "

# Row 47
$ws.Range("D47").Value = "'335.36"
$ws.Range("E47").Value = "'  -3.20%  "

# Row 48
$ws.Range("E48").Value = "'  -1.49%  "

# Row 49
$ws.Range("D49").Value = "'6.31"
$ws.Range("E49").Value = "'  +2.18%  "

# Row 50
$ws.Range("D50").Value = "'0.101"
$ws.Range("E50").Value = "'  -1.23%  "

# Row 51
$ws.Range("D51").Value = "'0.978"
$ws.Range("E51").Value = "'  -1.32%  "

